$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peer  and self assessment")

# Row 2 - "Self assesment" (online collaboration criterion, first block)
$ws.Range("B2").Value = "Good"
$ws.Range("C2").Value = "I have been active in the weekly online Discord meetings - both in `nterms of planning (Adding discussion points to the meeting agenda, taking responsibility for leading a group discussion or a walkthrough of conducted work, etc) `nI have been reading through literature and materials uploaded`nto GiT by other project group members, and actively been adding comments and suggestions. I have been uploading relevant materials to`nGiT as well.  `nI have been responding fairly quickly to messages, both private and group announcements, on the discord platform. "

# Row 22 - "Alex" (second block, international collaboration)
$ws.Range("B22").Value = "Excellent"
$ws.Range("C22").Value = "Active collaborator, motivated"

# Row 3 - "Alex" (first block, online collaboration)
$ws.Range("B3").Value = "Good"
$ws.Range("C3").Value = "Research, hardware setup "

# Row 15 - "Self assesment" (second block, international collaboration)
$ws.Range("B15").Value = "Good"
$ws.Range("C15").Value = "Active collaboration with Ahmet about Machine learning (Including sharing  research and literature), Active collaboration with Alex and Morcel about the PoA (Including giving inputs and suggestions for conducted work), Active communication and use of GiT and Discord, "

# Move the active selection to B16, matching the author's final cursor position
$null = $ws.Range("B16").Select()
